$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Pasta"

$ws.Range("C4").Value = "Watch movies"
$ws.Range("D4").Value = "Listen to music"
$ws.Range("E4").Value = "Go for walks"

$ws.Range("C6").Value = "The Secret"
